$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.066115702479338803
$ws.Range("B3").Value = 0.14049586776859499
$ws.Range("B4").Value = 0.23140495867768501
$ws.Range("B5").Value = 0.32231404958677601
